$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 219, pushing existing rows 219:289 down to 220:290
$ws.Rows.Item(219).Insert()

# Populate the new row 219 with the same "metadata" as the (now shifted) old
# row 219 (now at row 220), but with fresh data values reported in this commit.
$ws.Range("A219").Value = 10
$ws.Range("B219").Value = "Vega Modelo de Temuco"
$ws.Range("C219").Value = "La Araucanía"
$ws.Range("D219").Value = 44588
$ws.Range("E219").Value = 9
$ws.Range("F219").Value = "Fruta"
$ws.Range("G219").Value = 100108
$ws.Range("H219").Value = "Tropicales y subtropicales"
$ws.Range("I219").Value = 100108002
$ws.Range("J219").Value = "Mango"
$ws.Range("K219").Value = "Sin especificar"
$ws.Range("L219").Value = "Primera"
$ws.Range("M219").Value = 1200
$ws.Range("N219").Value = 8000
$ws.Range("O219").Value = 8000
$ws.Range("P219").Value = 8000
$ws.Range("Q219").Value = "$/bandeja 4 kilos"
$ws.Range("R219").Value = "Perú"
$ws.Range("S219").Value = 2000
$ws.Range("T219").Value = 4
